$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column to remain Text so numeric-looking strings
# (e.g. "0.4611", "11.95") are not auto-converted into Numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "28.856.00"
$ws.Range("D3").Value = "1.877.66"
$ws.Range("D5").Value = "324.39"
$ws.Range("D7").Value = "0.4611"
$ws.Range("D8").Value = "0.3871"
$ws.Range("D9").Value = "0.07845"
$ws.Range("D10").Value = "0.9833"
$ws.Range("D12").Value = "1.934.36"
$ws.Range("D13").Value = "6.990"
$ws.Range("D14").Value = "5.659"
$ws.Range("D15").Value = "0.06988"
$ws.Range("D16").Value = "88.16"
$ws.Range("D18").Value = "0.000009951"
$ws.Range("D19").Value = "16.91"
$ws.Range("D21").Value = "28.856.71"
$ws.Range("D22").Value = "5.256"
$ws.Range("D23").Value = "10.96"
$ws.Range("D25").Value = "156.44"
$ws.Range("D26").Value = "19.31"
$ws.Range("D27").Value = "5.958"
$ws.Range("D28").Value = "117.62"
$ws.Range("D29").Value = "1.906"
$ws.Range("D30").Value = "0.09353"
$ws.Range("D31").Value = "0.8997"
$ws.Range("D32").Value = "5.265"
$ws.Range("D34").Value = "3.249"
$ws.Range("D35").Value = "1.168"
$ws.Range("D37").Value = "0.02070"
$ws.Range("D39").Value = "7.628"
$ws.Range("D41").Value = "0.1769"
$ws.Range("D42").Value = "9.681"
$ws.Range("D43").Value = "11.95"
$ws.Range("D44").Value = "2.233"
$ws.Range("D45").Value = "0.5335"
$ws.Range("D47").Value = "1.838"
$ws.Range("D48").Value = "2.534"
$ws.Range("D49").Value = "112.40"
$ws.Range("D50").Value = "1.057"
$ws.Range("D51").Value = "70.95"

# Restore original (default/general) formatting now that the text values
# are safely stored, so no stray style attributes are left on the cells.
$priceRange.ClearFormats()

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("E3").Value = "  -1.95%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("E10").Value = "  -3.63%  "
$ws.Range("E11").Value = "  -2.77%  "
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("E19").Value = "  -2.85%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("E29").Value = "  -6.39%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  -4.75%  "
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("E35").Value = "  -0.97%  "
$ws.Range("E36").Value = "  -2.72%  "
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  -6.23%  "
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("E42").Value = "  -4.61%  "
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("E44").Value = "  -4.25%  "
$ws.Range("E45").Value = "  -2.92%  "
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("E47").Value = "  -3.89%  "
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("E49").Value = "  -1.08%  "
$ws.Range("E50").Value = "  -6.57%  "
$ws.Range("E51").Value = "  -1.19%  "
